$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60: politeness_score (B60) becomes a real number - was stored as text "3"
$ws.Range("B60").Value = 3

# New row 61: additional annotation entry for Ying Tang
$ws.Range("A61").Value = "Ying Tang"

# B61 keeps politeness_score stored as text "3" (matches the pre-edit B60 pattern)
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = "3"
$ws.Range("B61").Style = "Normal"

$ws.Range("C61").Value = " In any case,should be clarified."
$ws.Range("D61").Value = "SUG"
$ws.Range("E61").Value = "WRI"
$ws.Range("F61").Value = "1e0176d5-be35-49c3-adce-f7bfa3b6964a"
$ws.Range("G61").Value = "HksxTdiWz_annotated.xlsx"
$ws.Range("H61").Value = "In any case, this statement should be clarified."
